$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a roster in rows 2-10 (row 1 is the header, row 10 - the
# "Chiv" record - is unaffected). A new record ("Russ"/"Daren", with real
# scores instead of "NA") is inserted at row 2 and the rest of the roster
# (previously rows 2-8) shifts down to rows 3-9. The old duplicate "Russ"/
# "Daren" placeholder row (old row 9) is dropped since it is superseded by
# the new row 2.
#
# Values are written directly (bottom-up) instead of using Rows.Insert so
# that no formatting/style is copied from neighboring rows.

$data = @(
    @(18, "Singh",      "Ishdeep",  2,    2,    2),
    @(10, "Labiche",    "Yvan",     "NA", "NA", "NA"),
    @(11, "McConnell",  "Jenna",    "NA", "NA", "NA"),
    @(12, "Poll",       "Jennifer", "NA", "NA", "NA"),
    @(13, "Warmington", "Saundra",  "NA", "NA", "NA"),
    @(14, "East",       "Erica",    "NA", "NA", "NA"),
    @(15, "Buburuz",    "Jerry ",   "NA", "NA", "NA")
)

for ($i = $data.Length - 1; $i -ge 0; $i--) {
    $row = $i + 3
    $rec = $data[$i]
    $ws.Range("A$row").Value = $rec[0]
    $ws.Range("B$row").Value = $rec[1]
    $ws.Range("C$row").Value = $rec[2]
    $ws.Range("D$row").Value = $rec[3]
    $ws.Range("E$row").Value = $rec[4]
    $ws.Range("F$row").Value = $rec[5]
}

# New row 2: the "Russ"/"Daren" record with updated scores.
$ws.Range("A2").Value = 16
$ws.Range("B2").Value = "Russ"
$ws.Range("C2").Value = "Daren "
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 4
